# Apply the authored change set to SA_Info.xlsx:
#  - Append 8 new XPath-mapping rows (103-110) to the "XPath" sheet.
#  - Make "XPath" the active sheet / active tab (was "OrderInfo").
#  - Leave the scroll position near the newly-added rows with B112 selected.

$wb = $excel.ActiveWorkbook
$xpath = $wb.Worksheets.Item("XPath")

# New label/xpath pairs appended at the bottom of the XPath sheet.
$newRows = @(
    @("hide selectors", '//*[@id="root"]/div/div[4]/div[1]/div/div/div/div[2]/div[1]/div[1]/div[1]/div/div[2]/button/span'),
    @("cash back item from order detail tab", '//*[@id="root"]/div/div[4]/div[1]/div/div/div/div[2]/div[1]/div[3]/div[2]/div/div[2]/div/div/div[2]/div[2]/div/div/div/div[1]/div[1]/div[2]'),
    @("order detail tab ", '//*[@id="root"]/div/div[4]/div[1]/div/div/div/div[2]/div[1]/div[3]/div[2]/div/div[1]/ul/li[1]'),
    @("item_detail_tab", '//*[@id="root"]/div/div[4]/div[1]/div/div/div/div[2]/div[1]/div[3]/div[2]/div/div[1]/ul/li[1]/div'),
    @("logistic detail tab", '//*[@id="root"]/div/div[4]/div[1]/div/div/div/div[2]/div[1]/div[3]/div[2]/div/div[1]/ul/li[2]/div'),
    @("seller detail tab", '//*[@id="root"]/div/div[4]/div[1]/div/div/div/div[2]/div[1]/div[3]/div[2]/div/div[1]/ul/li[3]/div/div/img'),
    @("offer detail tab", '//*[@id="root"]/div/div[4]/div[1]/div/div/div/div[2]/div[1]/div[3]/div[2]/div/div[1]/ul/li[4]/div/span'),
    @("item detail all DT's", '//*[@id="root"]/div/div[4]/div[1]/div/div/div/div[2]/div[1]/div[3]/div[2]/div/div[2]/div/div/div[1]/div[2]')
)

$startRow = 103
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $xpath.Cells.Item($r, 1).Value = $newRows[$i][0]
    $xpath.Cells.Item($r, 2).Value = $newRows[$i][1]
}

# Fix the curly apostrophe in the last label (PowerShell source above used a
# plain apostrophe to avoid quoting issues).
$xpath.Cells.Item(110, 1).Value = "item detail all DT" + [char]0x2019 + "s"

# Make XPath the active sheet/tab and position the view near the new rows.
$xpath.Activate()
$excel.Goto($xpath.Range("B95"), $true)
$xpath.Range("B112").Select()
